# Applies the "feat: add 2022-Q3 data" change:
#  1. Inserts a new "2022-Q3" quarterly data row at the top of the "总计"
#     summary sheet (pushing the other quarters down, and appending the
#     former last row, 2021-Q1, as a brand-new trailing row).
#  2. Inserts a brand-new "2022-Q3" worksheet (a duplicate of "2022-Q2"'s
#     layout/format) right before "2022-Q2", populated with the Q3 fund
#     holding data.

function Set-TextValue($cell, $val) {
    # Force the cell to be stored as literal text (not auto-coerced to a
    # number), matching the source file's inlineStr cells for numeric-looking
    # strings like fund codes ("012703") and percentages ("4.94"), while
    # leaving the cell's style untouched (no lingering number-format style).
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert the 2022-Q3 summary row, shifting the rest down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Remember the current last data row (2021-Q1, row 7) before it gets
# overwritten by the downward shift.
$lastB = $total.Cells.Item(7, 2).Value2
$lastC = $total.Cells.Item(7, 3).Value2
$lastD = $total.Cells.Item(7, 4).Value2

# Create the new trailing row 8, copying column A's formatting (s="2")
# from the row directly above it.
$total.Cells.Item(8, 1).Value = 6
$total.Cells.Item(7, 1).Copy()
$total.Cells.Item(8, 1).PasteSpecial(-4122)
$total.Cells.Item(8, 2).Value = $lastB
$total.Cells.Item(8, 3).Value = $lastC
$total.Cells.Item(8, 4).Value = $lastD

# Shift existing rows 2-6 down into rows 3-7.
for ($r = 7; $r -ge 3; $r--) {
    $src = $r - 1
    $total.Cells.Item($r, 2).Value = $total.Cells.Item($src, 2).Value2
    $total.Cells.Item($r, 3).Value = $total.Cells.Item($src, 3).Value2
    $total.Cells.Item($r, 4).Value = $total.Cells.Item($src, 4).Value2
}

# Write the new 2022-Q3 row into row 2 (column A / style already correct).
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 10
$total.Cells.Item(2, 4).Value = 0.41

# ---------------------------------------------------------------------
# 2. Add the new "2022-Q3" worksheet before "2022-Q2", cloning its
#    layout/formatting and filling in the Q3 fund-holding data.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The source sheet (2022-Q2) has 16 data rows; 2022-Q3 only needs 10, so
# drop the extra trailing rows (12-17).
$q3.Rows("12:17").Delete()

$rows = @(
    @(0, "012703", "华夏核心成长混合A", "4.94", "93.34", "4.55", "0.2248", 9),
    @(1, "011097", "达诚宜创精选混合A", "0.70", "73.65", "5.63", "0.0394", 3),
    @(2, "010301", "达诚成长先锋混合A", "0.52", "74.13", "5.69", "0.0296", 3),
    @(3, "012710", "华夏核心成长混合C", "0.61", "93.34", "4.55", "0.0278", 9),
    @(4, "010808", "达诚策略先锋混合A", "0.33", "75.14", "5.74", "0.0189", 3),
    @(5, "010809", "达诚策略先锋混合C", "0.32", "75.14", "5.74", "0.0184", 3),
    @(6, "010302", "达诚成长先锋混合C", "0.31", "74.13", "5.69", "0.0176", 3),
    @(7, "011031", "达诚价值先锋灵活配置混合C", "0.23", "75.36", "5.49", "0.0126", 3),
    @(8, "011030", "达诚价值先锋灵活配置混合A", "0.21", "75.36", "5.49", "0.0115", 3),
    @(9, "011098", "达诚宜创精选混合C", "0.18", "73.65", "5.63", "0.0101", 3)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    Set-TextValue $q3.Cells.Item($r, 2) $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    Set-TextValue $q3.Cells.Item($r, 4) $row[3]
    Set-TextValue $q3.Cells.Item($r, 5) $row[4]
    Set-TextValue $q3.Cells.Item($r, 6) $row[5]
    Set-TextValue $q3.Cells.Item($r, 7) $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
